$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers. The order in which new header strings are first written controls
# the order they land in the shared-string table, so write "import " and
# "total" before "swiss license".
$ws.Range("A1").Value = "team"
$ws.Range("C1").Value = "import "
$ws.Range("D1").Value = "total"
$ws.Range("B1").Value = "swiss license"

# Per-team swiss license / import / total numbers.
$swiss = @(21,23,20,22,21,22,21,21,19,22,20,21,25,22)
$import = @(7,4,4,7,7,5,4,5,6,4,4,5,2,5)
$total = @(28,27,24,29,28,27,25,26,25,26,24,26,27,27)

for ($i = 0; $i -lt $swiss.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $swiss[$i]
    $ws.Cells.Item($row, 3).Value = $import[$i]
    $ws.Cells.Item($row, 4).Value = $total[$i]
}

# Totals row
$ws.Range("B16").Value = 300
$ws.Range("C16").Formula = "=SUM(C2:C15)"
$ws.Range("D16").Formula = "=SUM(D2:D15)"

# Column widths: widen the custom width (originally applied only to A) to
# also cover the new B and C columns, keeping A's own width untouched.
$ws.Columns("B:C").ColumnWidth = 23.35

# Selection, to mirror the saved view state
$ws.Range("C9").Select() | Out-Null
